# Update header row labels so the first row can be auto-detected as a
# header when the data is loaded into Power BI.
# Sheets 1,2,3,5 (Potencia Acumulada, Geracao Periodo Medio, Atendimento a
# Ponta, Emissoes Totais) use year columns -> prefix with "Ano ".
# Sheet 4 (Potencia Incremental) uses year/interval columns -> prefix with
# "Intervalo ".
# Sheet 6 (Custo Total) only has a single year column (B1) -> prefix with
# "Ano ".

$wb = $excel.ActiveWorkbook

$anoSheets = @(1, 2, 3, 5)
foreach ($idx in $anoSheets) {
    $ws = $wb.Worksheets.Item($idx)
    $ws.Range("B1").Value = "Ano " + $ws.Range("B1").Text
    $ws.Range("C1").Value = "Ano " + $ws.Range("C1").Text
    $ws.Range("D1").Value = "Ano " + $ws.Range("D1").Text
    $ws.Range("E1").Value = "Ano " + $ws.Range("E1").Text
}

$wsIntervalo = $wb.Worksheets.Item(4)
$wsIntervalo.Range("B1").Value = "Intervalo " + $wsIntervalo.Range("B1").Text
$wsIntervalo.Range("C1").Value = "Intervalo " + $wsIntervalo.Range("C1").Text
$wsIntervalo.Range("D1").Value = "Intervalo " + $wsIntervalo.Range("D1").Text
$wsIntervalo.Range("E1").Value = "Intervalo " + $wsIntervalo.Range("E1").Text

$wsCusto = $wb.Worksheets.Item(6)
$wsCusto.Range("B1").Value = "Ano " + $wsCusto.Range("B1").Text
